$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 96, shifting rows 96:126 down to 97:127
$ws.Rows.Item(96).Insert()

# Populate the new row 96 with the latest weekly price entry
$ws.Cells.Item(96, 1).Value = 10
$ws.Cells.Item(96, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(96, 3).Value = "La Araucanía"
$ws.Cells.Item(96, 4).Value = 44785
$ws.Cells.Item(96, 5).Value = 9
$ws.Cells.Item(96, 6).Value = "Fruta"
$ws.Cells.Item(96, 7).Value = 100107
$ws.Cells.Item(96, 8).Value = "Otros"
$ws.Cells.Item(96, 9).Value = 100107002
$ws.Cells.Item(96, 10).Value = "Chirimoya"
$ws.Cells.Item(96, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(96, 12).Value = "Primera"
$ws.Cells.Item(96, 13).Value = 25
$ws.Cells.Item(96, 14).Value = 4000
$ws.Cells.Item(96, 15).Value = 4000
$ws.Cells.Item(96, 16).Value = 4000
$ws.Cells.Item(96, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(96, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(96, 19).Value = 4000
$ws.Cells.Item(96, 20).Value = 1
